$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) "Right click: Set P2/Offset" paragraph loses the _GoBack bookmark
#    (it gets moved to the new "Left/Right arrow" paragraph below).
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) "Enter: Accept current trial (closing figure window will also do
#    this)" -> multi-run "Enter: Accept all trials (closing figure
#    window will also do this)"
# ------------------------------------------------------------------
$pEnter = $d.Paragraphs.Item(4)
$xmlEnter = '<w:p ' + $wNs + '>' +
  '<w:r><w:t xml:space="preserve">Enter: Accept </w:t></w:r>' +
  '<w:r><w:t>all</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> trial</w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> (closing figure window will also do this)</w:t></w:r>' +
  '</w:p>'
$pEnter.Range.InsertXML($xmlEnter)

# ------------------------------------------------------------------
# 3) Remove the "Ctrl+Enter: Accept current and all remaining trials"
#    paragraph entirely.
# ------------------------------------------------------------------
$pCtrlEnter = $d.Paragraphs.Item(5)
$pCtrlEnter.Range.Delete()

# ------------------------------------------------------------------
# 4) "Delete or Backspace: Reject current trial" -> two-run
#    "Delete or Backspace: Remove peaks from current trial"
# ------------------------------------------------------------------
$pDelete = $d.Paragraphs.Item(5)
$xmlDelete = '<w:p ' + $wNs + '>' +
  '<w:r><w:t xml:space="preserve">Delete or Backspace: </w:t></w:r>' +
  '<w:r><w:t>Remove peaks from current trial</w:t></w:r>' +
  '</w:p>'
$pDelete.Range.InsertXML($xmlDelete)

# ------------------------------------------------------------------
# 5) Insert the new "Left/Right arrow: Cycle backward/forward through
#    trials" paragraph right after, carrying the relocated _GoBack
#    bookmark at its end (after the run, like the original placement).
# ------------------------------------------------------------------
$pDelete = $d.Paragraphs.Item(5)
$pDelete.Range.InsertParagraphAfter()

$pArrow = $d.Paragraphs.Item(6)
$xmlArrow = '<w:p ' + $wNs + '>' +
  '<w:r><w:t>Left/Right arrow: Cycle backward/forward through trials</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$pArrow.Range.InsertXML($xmlArrow)
